$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 - copy the existing header formatting
# from G1 (bold font, border, centered alignment) so the new cell reuses
# the same style, then set its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column with 0 for each data row (2-6), matching
# the plain (unstyled) numeric cells already used in columns B-G.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
